$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("RUIZ CHIROQUE CLAUDIA JUDITH", 174),
    @("FABIANA REBECA ARRUNATEGUI SILUPU", 161),
    @("BANCAYAN FIESTA DILVER HUMBERTO", 158),
    @("LLENQUE ANTON HELEN JOHANA", 157),
    @("TEMOCHE ECHE URSULA YESSENIA", 153),
    @("GONZALES FIESTAS MARIA MARIBEL", 152),
    @("ANTON INGA FATIMA DEL ROSARIO", 150),
    @("BAUTISTA CHAVESTA ERICKA MEDALIT", 142),
    @("VELASCO PEÑA KAREN ARELLYS", 136),
    @("MONDRAGON NONAJULCA MARISOL", 124),
    @("PINTADO CHASQUERO ESTEFANY", 117),
    @("FLORES SILUPU MARY CARMEN", 116),
    @("ORDINOLA JIBAJA JOSE ALBERTO", 113),
    @("HERNANDEZ CARNERO ARTURO SEBASTIAN", 109),
    @("MORENO YANAYACO NAYLA GUADALUPE", 89),
    @("CASTRO ESTRADA CINTHIA PATRICIA", 84),
    @("PINTADO BENITES CRISTOBAL RODRIGO", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
